# Auto-generated Excel COM-interop edit script
# Applies the row-content permutation + field edits for rows 4-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 100891113
$ws.Range("B4").Value = 56395
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("M4").Value = "färska spår"
$ws.Range("Q4").Value = 806814.4326894956
$ws.Range("R4").Value = 7175557.99483321
$ws.Range("A5").Value = 100945295
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("P5").Value = "Burvik, Vb"
$ws.Range("Q5").Value = 806845.0991897933
$ws.Range("R5").Value = 7175490.583056876
$ws.Range("S5").Value = 25
$ws.Range("AW5").Value = "Patrik Nygren"
$ws.Range("AX5").Value = "Patrik Nygren"
$ws.Range("A6").Value = 100945131
$ws.Range("B6").Value = 81236
$ws.Range("E6").Value = 1312
$ws.Range("F6").Value = "Gammelgransskål"
$ws.Range("G6").Value = "Pseudographis pinicola"
$ws.Range("H6").Value = "(Nyl.) Rehm"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("P6").Value = "Burvik, Vb"
$ws.Range("Q6").Value = 806795.2342230955
$ws.Range("R6").Value = 7175644.739797495
$ws.Range("S6").Value = 25
$ws.Range("AW6").Value = "Patrik Nygren"
$ws.Range("AX6").Value = "Patrik Nygren"
$ws.Range("A7").Value = 101279776
$ws.Range("B7").Value = 96334
$ws.Range("D7").Value = "VU"
$ws.Range("E7").Value = 220787
$ws.Range("F7").Value = "Knärot"
$ws.Range("G7").Value = "Goodyera repens"
$ws.Range("H7").Value = "(L.) R. Br."
$ws.Range("I7").Value = "1"
$ws.Range("J7").Value = "plantor/tuvor"
$ws.Range("P7").Value = "Burvik, Burvik, Vb"
$ws.Range("Q7").Value = 806838.9976605003
$ws.Range("R7").Value = 7175499.482516243
$ws.Range("Y7").Value = "2022-05-24"
$ws.Range("Z7").Value = "10:00"
$ws.Range("AA7").Value = "2022-05-24"
$ws.Range("AB7").Value = "12:00"
$ws.Range("AW7").Value = "Jon Andersson"
$ws.Range("AX7").Value = "Jon Andersson, Patrik Nygren"
$ws.Range("A8").Value = 100890468
$ws.Range("B8").Value = 89356
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 5447
$ws.Range("F8").Value = "Vedticka"
$ws.Range("G8").Value = "Fuscoporia viticola"
$ws.Range("H8").Value = "(Schwein.) Murrill"
$ws.Range("P8").Value = "Karlgrundsmoren, Vb"
$ws.Range("Q8").Value = 806772.363243226
$ws.Range("R8").Value = 7175558.485492887
$ws.Range("S8").Value = 10
$ws.Range("AJ8").Value = "gran"
$ws.Range("AK8").Value = "Picea abies"
$ws.Range("AO8").Value = "Picea abies"
$ws.Range("AW8").Value = "Roger Olofsson"
$ws.Range("AX8").Value = "Roger Olofsson"
$ws.Range("A9").Value = 100887736
$ws.Range("B9").Value = 89392
$ws.Range("E9").Value = 1202
$ws.Range("F9").Value = "Ullticka"
$ws.Range("G9").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H9").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P9").Value = "Karlgrundsmoren, Vb"
$ws.Range("Q9").Value = 806712.6351854501
$ws.Range("R9").Value = 7175617.756130967
$ws.Range("S9").Value = 10
$ws.Range("AC9").Value = "Torkad"
$ws.Range("AJ9").Value = "gran"
$ws.Range("AK9").Value = "Picea abies"
$ws.Range("AO9").Value = "Picea abies"
$ws.Range("AW9").Value = "Roger Olofsson"
$ws.Range("AX9").Value = "Roger Olofsson"
$ws.Range("A10").Value = 100890647
$ws.Range("B10").Value = 89356
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 5447
$ws.Range("F10").Value = "Vedticka"
$ws.Range("G10").Value = "Fuscoporia viticola"
$ws.Range("H10").Value = "(Schwein.) Murrill"
$ws.Range("P10").Value = "Karlgrundsmoren, Vb"
$ws.Range("Q10").Value = 806785.6405356181
$ws.Range("R10").Value = 7175538.632476954
$ws.Range("S10").Value = 10
$ws.Range("AJ10").Value = "gran"
$ws.Range("AK10").Value = "Picea abies"
$ws.Range("AO10").Value = "Picea abies"
$ws.Range("AW10").Value = "Roger Olofsson"
$ws.Range("AX10").Value = "Roger Olofsson"
$ws.Range("A11").Value = 100945031
$ws.Range("Q11").Value = 806750.2499365495
$ws.Range("R11").Value = 7175580.037869117
$ws.Range("A12").Value = 100945368
$ws.Range("B12").Value = 89388
$ws.Range("E12").Value = 1108
$ws.Range("F12").Value = "Harticka"
$ws.Range("G12").Value = "Pelloporus leporinus"
$ws.Range("H12").Value = "(Fr.) Krieglst."
$ws.Range("Q12").Value = 806874.6453256523
$ws.Range("R12").Value = 7175459.833355412
$ws.Range("A13").Value = 100945088
$ws.Range("B13").Value = 89392
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q13").Value = 806765.2456494838
$ws.Range("R13").Value = 7175590.212229759
$ws.Range("A14").Value = 100945420
$ws.Range("B14").Value = 81236
$ws.Range("E14").Value = 1312
$ws.Range("F14").Value = "Gammelgransskål"
$ws.Range("G14").Value = "Pseudographis pinicola"
$ws.Range("H14").Value = "(Nyl.) Rehm"
$ws.Range("Q14").Value = 806950.0357482962
$ws.Range("R14").Value = 7175429.842739396
$ws.Range("A15").Value = 100945401
$ws.Range("B15").Value = 89392
$ws.Range("E15").Value = 1202
$ws.Range("F15").Value = "Ullticka"
$ws.Range("G15").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H15").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("J15").Value = "fruktkroppar"
$ws.Range("Q15").Value = 806914.5347780141
$ws.Range("R15").Value = 7175455.228436214
$ws.Range("A16").Value = 100945319
$ws.Range("B16").Value = 77506
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = "Garnlav"
$ws.Range("G16").Value = "Alectoria sarmentosa"
$ws.Range("H16").Value = "(Ach.) Ach."
$ws.Range("J16").Value = "bålar"
$ws.Range("Q16").Value = 806849.8072187628
$ws.Range("R16").Value = 7175465.533442252
$ws.Range("A17").Value = 100945109
$ws.Range("B17").Value = 89356
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 5447
$ws.Range("F17").Value = "Vedticka"
$ws.Range("G17").Value = "Fuscoporia viticola"
$ws.Range("H17").Value = "(Schwein.) Murrill"
$ws.Range("J17").Value = "fruktkroppar"
$ws.Range("P17").Value = "Burvik, Vb"
$ws.Range("Q17").Value = 806775.5716160003
$ws.Range("R17").Value = 7175629.332793099
$ws.Range("S17").Value = 25
$ws.Range("Y17").Value = "2022-05-15"
$ws.Range("Z17").Value = "00:00"
$ws.Range("AA17").Value = "2022-05-15"
$ws.Range("AB17").Value = "00:00"
$ws.Range("AW17").Value = "Patrik Nygren"
$ws.Range("AX17").Value = "Patrik Nygren"

# Clear cells that should become empty/absent
$ws.Range("J4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("AF4").ClearContents()
$ws.Range("AC5").ClearContents()
$ws.Range("AJ5").ClearContents()
$ws.Range("AK5").ClearContents()
$ws.Range("AO5").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()
$ws.Range("AF6").ClearContents()
$ws.Range("AJ6").ClearContents()
$ws.Range("AK6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("AJ7").ClearContents()
$ws.Range("AK7").ClearContents()
$ws.Range("AO7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("I17").ClearContents()
$ws.Range("L17").ClearContents()
